$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69; this shifts old rows 69-124 down to 70-125,
# matching the new sheet dimension A1:T125.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the latest weekly price entry.
$ws.Cells.Item(69,1).Value = 7
$ws.Cells.Item(69,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69,3).Value = "Ñuble"
$ws.Cells.Item(69,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(69,4).Value = 44902
$ws.Cells.Item(69,5).Value = 16
$ws.Cells.Item(69,6).Value = "Fruta"
$ws.Cells.Item(69,7).Value = 100103
$ws.Cells.Item(69,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(69,9).Value = 100103001
$ws.Cells.Item(69,10).Value = "Cereza"
$ws.Cells.Item(69,11).Value = "Santina"
$ws.Cells.Item(69,12).Value = "Primera"
$ws.Cells.Item(69,13).Value = 120
$ws.Cells.Item(69,14).Value = 6000
$ws.Cells.Item(69,15).Value = 7000
$ws.Cells.Item(69,16).Value = 6500
$ws.Cells.Item(69,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(69,18).Value = "Provincia de Curicó"
$ws.Cells.Item(69,19).Value = 650
$ws.Cells.Item(69,20).Value = 10
